# TestRunner.xlsx update
# 1. Rename the "NumberOfWindows" config property label to "NumberOfBrowsers".
# 2. Make the "Config" sheet the active/selected sheet (was "Test Cases"),
#    and move its selection from B6 to A7.

$wb = $excel.ActiveWorkbook

$configSheet = $wb.Worksheets.Item("Config")

# Rename the config property label in column A, row 4.
$configSheet.Range("A4").Value = "NumberOfBrowsers"

# Switch the active tab to "Config" and move its selection to A7.
$configSheet.Activate()
$configSheet.Range("A7").Select()

# Leave "Test Cases" sheet's own selection untouched (still D10).
